$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 20 data - copy formatting from row above (A19) then set the new value
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A20").Value = 45986

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 0.8976398032236155
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 0.4275768375374467
